$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-17 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-18 Monday", 2)

$d.Content.Find.Execute("77×50=3850", $true, $false, $false, $false, $false, $true, 1, $false, "84×24=2016", 2)
$d.Content.Find.Execute("37×16=592", $true, $false, $false, $false, $false, $true, 1, $false, "41×75=3075", 2)
$d.Content.Find.Execute("30×23=690", $true, $false, $false, $false, $false, $true, 1, $false, "59×69=4071", 2)
$d.Content.Find.Execute("89×96=8544", $true, $false, $false, $false, $false, $true, 1, $false, "18×47=846", 2)
$d.Content.Find.Execute("55×17=935", $true, $false, $false, $false, $false, $true, 1, $false, "12×57=684", 2)
$d.Content.Find.Execute("62×41=2542", $true, $false, $false, $false, $false, $true, 1, $false, "53×34=1802", 2)
$d.Content.Find.Execute("89×83=7387", $true, $false, $false, $false, $false, $true, 1, $false, "17×25=425", 2)
$d.Content.Find.Execute("26×39=1014", $true, $false, $false, $false, $false, $true, 1, $false, "38×55=2090", 2)
$d.Content.Find.Execute("30×12=360", $true, $false, $false, $false, $false, $true, 1, $false, "12×20=240", 2)
$d.Content.Find.Execute("12×71=852", $true, $false, $false, $false, $false, $true, 1, $false, "46×30=1380", 2)
$d.Content.Find.Execute("20×89=1780", $true, $false, $false, $false, $false, $true, 1, $false, "45×20=900", 2)
$d.Content.Find.Execute("25×60=1500", $true, $false, $false, $false, $false, $true, 1, $false, "28×33=924", 2)
$d.Content.Find.Execute("43×20=860", $true, $false, $false, $false, $false, $true, 1, $false, "28×61=1708", 2)
$d.Content.Find.Execute("13×94=1222", $true, $false, $false, $false, $false, $true, 1, $false, "25×81=2025", 2)
$d.Content.Find.Execute("97×49=4753", $true, $false, $false, $false, $false, $true, 1, $false, "45×47=2115", 2)
$d.Content.Find.Execute("41×70=2870", $true, $false, $false, $false, $false, $true, 1, $false, "50×76=3800", 2)
$d.Content.Find.Execute("81×99=8019", $true, $false, $false, $false, $false, $true, 1, $false, "27×16=432", 2)
$d.Content.Find.Execute("14×26=364", $true, $false, $false, $false, $false, $true, 1, $false, "29×94=2726", 2)
$d.Content.Find.Execute("21×45=945", $true, $false, $false, $false, $false, $true, 1, $false, "70×22=1540", 2)
$d.Content.Find.Execute("40×97=3880", $true, $false, $false, $false, $false, $true, 1, $false, "29×54=1566", 2)
$d.Content.Find.Execute("35×52=1820", $true, $false, $false, $false, $false, $true, 1, $false, "27×33=891", 2)
$d.Content.Find.Execute("47×85=3995", $true, $false, $false, $false, $false, $true, 1, $false, "88×25=2200", 2)
$d.Content.Find.Execute("79×55=4345", $true, $false, $false, $false, $false, $true, 1, $false, "29×79=2291", 2)
$d.Content.Find.Execute("29×42=1218", $true, $false, $false, $false, $false, $true, 1, $false, "61×93=5673", 2)
$d.Content.Find.Execute("91×32=2912", $true, $false, $false, $false, $false, $true, 1, $false, "39×82=3198", 2)
